$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "seatsNumber"
$ws.Range("B8").Value = 5

$ws.Range("A9").Value = "fuelType"
$ws.Range("B9").Value = "Diesel;Petrol"

$ws.Range("A10").Value = "powerNumber"
$ws.Range("B10").Value = 131

$ws.Range("A11").Value = "cubicCapacity"
$ws.Range("B11").Value = 3000

$ws.Range("A12").Value = "transmissionType"
$ws.Range("B12").Value = "Manual Gearbox"

$ws.Range("A13").Value = "emissionSticker"
$ws.Range("B13").Value = 4

$ws.Range("B14").Select()
